$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.916.39'
$ws.Range("E2").Value = '  -2.38%  '
$ws.Range("D3").Value = '2.579.92'
$ws.Range("E3").Value = '  -3.90%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '551.33'
$ws.Range("E5").Value = '  -0.54%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '154.66'
$ws.Range("E6").Value = '  -2.63%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.593'
$ws.Range("E8").Value = '  +1.17%  '
$ws.Range("E9").Value = '  -1.55%  '
$ws.Range("E10").Value = '  -0.80%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.50'
$ws.Range("E11").Value = '  +3.39%  '
$ws.Range("E12").Value = '  -0.81%  '
$ws.Range("D13").Value = '3.034.94'
$ws.Range("E13").Value = '  -3.89%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '25.45'
$ws.Range("E14").Value = '  -3.74%  '
$ws.Range("D15").Value = '61.843.44'
$ws.Range("E15").Value = '  -2.23%  '
$ws.Range("E16").Value = '  -0.89%  '
$ws.Range("D17").Value = '2.584.98'
$ws.Range("E17").Value = '  -3.76%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '11.59'
$ws.Range("E18").Value = '  -3.53%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.53'
$ws.Range("E19").Value = '  -0.94%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '338.56'
$ws.Range("E20").Value = '  -2.09%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.01'
$ws.Range("E21").Value = '  -4.78%  '
$ws.Range("E22").Value = '  +0.24%  '
$ws.Range("E23").Value = '  -3.36%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '63.54'
$ws.Range("E24").Value = '  -0.42%  '
$ws.Range("E25").Value = '  -1.11%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.999'
$ws.Range("E26").Value = '  +0.03%  '
$ws.Range("E27").Value = '  -0.75%  '
$ws.Range("E28").Value = '  +2.98%  '
$ws.Range("E29").Value = '  -3.17%  '
$ws.Range("E30").Value = '  -1.23%  '
$ws.Range("E31").Value = '  -2.85%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '162.75'
$ws.Range("E32").Value = '  -1.90%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.87'
$ws.Range("E33").Value = '  +0.99%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '19.13'
$ws.Range("E35").Value = '  -2.14%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.40'
$ws.Range("E36").Value = '  -2.28%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.78'
$ws.Range("E37").Value = '  -0.46%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.03'
$ws.Range("E38").Value = '  -1.54%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '326.41'
$ws.Range("E39").Value = '  -4.97%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.899'
$ws.Range("E40").Value = '  -4.84%  '
$ws.Range("E41").Value = '  -0.14%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '37.48'
$ws.Range("E42").Value = '  -1.50%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '20.59'
$ws.Range("E43").Value = '  -1.21%  '
$ws.Range("E44").Value = '  +0.00%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.606'
$ws.Range("E46").Value = '  -1.08%  '
$ws.Range("E47").Value = '  -3.29%  '
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '19.51'
$ws.Range("E48").Value = '  -4.04%  '
$ws.Range("B49").Value = 'Stellar'
$ws.Range("C49").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0964'
$ws.Range("E49").Value = '  -0.86%  '
$ws.Range("E50").Value = '  -1.39%  '
$ws.Range("D51").Value = '2.049.10'
$ws.Range("E51").Value = '  -2.01%  '
